# "Zeitblatt auf stand gebracht" - update the October ("Oktober") timesheet
# with a new day's entries (row 8) and extend detail notes for rows 5-7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oktober")

# New day entry (05.10.2016, row 8): time range, description, hours worked
$ws.Range("D8").Value = "13:15 - 16:30, 17:00 - 19:30"
$ws.Range("C8").Value = "Teammeeting, Fertigstellen von Präsi"
$ws.Range("D7").Value = "19:00 - 20:00"

# Extra "Details" time-range notes for the days already present (rows 5-6)
$ws.Range("D5").Value = "17:00 - 21:00"
$ws.Range("D6").Value = "21:00  - 24:00"

# Hours worked for the new day
$ws.Range("B8").Value = 5.75

# Keep the selection where it was left in the saved file
$ws.Range("D7").Select()
